# Commit: "prepare for human readabily"
#
# The "Reference" column (column A) of the Revelation sheet was generated
# with a stray "16" suffix appended to many verse references
# (e.g. "Revelation 1:116" instead of "Revelation 1:1"). This script
# strips that trailing "16" from every affected cell in column A so the
# references read correctly, leaving every other cell (including the
# "Text" column B, and rows whose reference never had the stray suffix)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (falls back to the sheet's used range).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $text = $cell.Text

    if ($text -ne $null -and $text.EndsWith("16")) {
        $fixed = $text.Substring(0, $text.Length - 2)
        $cell.Value = $fixed
    }
}
